$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update the product record (Duet South Sea Pearl Earrings -> Freshwater Pearl and Diamond Necklace) ---

# Reference Number
$ws.Range("B2").Value = "Q2091NLET-7F"
# Collection
$ws.Range("C2").Value = "Trend"
# Product Name
$ws.Range("D2").Value = "Freshwater Pearl and Diamond Necklace"
$ws.Range("D2").HorizontalAlignment = -4131
# Description
$ws.Range("E2").Value = "Putting style at the forefront of each design, our Trend collection is perfect for the modern woman. This 18ct white gold necklace features high quality 4-6.5mm Freshwater pearls, accentuated by scintillating diamonds. A modern 'Y'-shaped silhouette is both elegant and contemporary, and will elevate any outfit. "
# Price (now a real number instead of "Call For Price")
$ws.Range("F2").Value = 1560
# Quantity
$ws.Range("H2").Value = 0
# Category
$ws.Range("K2").Value = "Necklace"
# Sub-Category (new, stays blank but becomes a real formatted cell)
$ws.Range("L2").NumberFormat = "General"
# Stone
$ws.Range("N2").Value = "Pearl - Diamond"
# Total Diamond Weight (new)
$ws.Range("P2").Value = "0.188ct"
# Also Available In (new)
$ws.Range("Q2").Value = "18K Yellow Gold"
# Pearl Type
$ws.Range("T2").Value = "Freshwater"
# Pearl Size
$ws.Range("U2").Value = "4 - 6.5 mm"
# Picture 1 (new)
$ws.Range("V2").Value = "Q2091NLET-7F_1000.png"
# Picture 2
$ws.Range("W2").Value = "Q2091NLET-7F CLASP_1000.png"
# Picture 3 (new)
$ws.Range("X2").Value = "Q2091NLET-7f M_1000.png"

# --- Expand the sheet's used range / dimension to A1:X1003 ---
$ws.Range("X1003").NumberFormat = "General"

# --- Update the active selection to B2 ---
$ws.Range("B2").Select()
